$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '37.656.68'
Set-TextValue $ws.Range("E2") '  +6.48%  '

Set-TextValue $ws.Range("D3") '2.044.34'
Set-TextValue $ws.Range("E3") '  +3.61%  '

Set-TextValue $ws.Range("D5") '251.17'
Set-TextValue $ws.Range("E5") '  +4.62%  '

Set-TextValue $ws.Range("E6") '  +2.84%  '

Set-TextValue $ws.Range("D7") '65.70'
Set-TextValue $ws.Range("E7") '  +17.45%  '

Set-TextValue $ws.Range("E8") '  -0.01%  '

Set-TextValue $ws.Range("D9") '59.62'
Set-TextValue $ws.Range("E9") '  +0.30%  '

Set-TextValue $ws.Range("E10") '  +5.64%  '

Set-TextValue $ws.Range("E11") '  +4.37%  '

Set-TextValue $ws.Range("E12") '  +0.97%  '

Set-TextValue $ws.Range("D13") '0.905'
Set-TextValue $ws.Range("E13") '  +1.80%  '

Set-TextValue $ws.Range("D14") '15.10'
Set-TextValue $ws.Range("E14") '  +6.84%  '

Set-TextValue $ws.Range("D15") '2.342.37'
Set-TextValue $ws.Range("E15") '  +3.40%  '

Set-TextValue $ws.Range("D16") '5.58'
Set-TextValue $ws.Range("E16") '  +6.95%  '

Set-TextValue $ws.Range("D17") '20.75'
Set-TextValue $ws.Range("E17") '  +21.61%  '

Set-TextValue $ws.Range("D18") '2.067.10'
Set-TextValue $ws.Range("E18") '  +4.12%  '

Set-TextValue $ws.Range("D19") '37.539.07'
Set-TextValue $ws.Range("E19") '  +6.61%  '

Set-TextValue $ws.Range("D20") '73.31'
Set-TextValue $ws.Range("E20") '  +5.07%  '

Set-TextValue $ws.Range("D21") '0.0₃0875'
Set-TextValue $ws.Range("E21") '  +5.21%  '

Set-TextValue $ws.Range("D22") '5.35'
Set-TextValue $ws.Range("E22") '  +6.76%  '

Set-TextValue $ws.Range("D23") '237.41'
Set-TextValue $ws.Range("E23") '  +2.63%  '

Set-TextValue $ws.Range("E24") '  +19.65%  '

Set-TextValue $ws.Range("E25") '  -0.03%  '

Set-TextValue $ws.Range("D26") '2.38'
Set-TextValue $ws.Range("E26") '  +5.87%  '

Set-TextValue $ws.Range("D27") '9.59'
Set-TextValue $ws.Range("E27") '  +6.45%  '

Set-TextValue $ws.Range("D28") '164.83'
Set-TextValue $ws.Range("E28") '  +1.13%  '

Set-TextValue $ws.Range("D29") '19.89'
Set-TextValue $ws.Range("E29") '  +2.91%  '

Set-TextValue $ws.Range("D30") '0.121'
Set-TextValue $ws.Range("E30") '  +3.28%  '

Set-TextValue $ws.Range("D31") '5.20'
Set-TextValue $ws.Range("E31") '  +9.69%  '

Set-TextValue $ws.Range("D32") '1.21'
Set-TextValue $ws.Range("E32") '  +8.26%  '

Set-TextValue $ws.Range("D33") '0.112'
Set-TextValue $ws.Range("E33") '  +25.32%  '

Set-TextValue $ws.Range("D34") '4.74'
Set-TextValue $ws.Range("E34") '  +11.91%  '

Set-TextValue $ws.Range("D35") '0.0612'
Set-TextValue $ws.Range("E35") '  +5.49%  '

Set-TextValue $ws.Range("D36") '2.45'
Set-TextValue $ws.Range("E36") '  +8.46%  '

Set-TextValue $ws.Range("E37") '  -0.10%  '

Set-TextValue $ws.Range("B38") 'WEMIXToken'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D38") '1.83'
Set-TextValue $ws.Range("E38") '  +2.45%  '

Set-TextValue $ws.Range("B39") 'THORChain'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D39") '6.08'
Set-TextValue $ws.Range("E39") '  +26.14%  '

Set-TextValue $ws.Range("E40") '  +16.98%  '

Set-TextValue $ws.Range("E41") '  +4.28%  '

Set-TextValue $ws.Range("E42") '  +22.98%  '

Set-TextValue $ws.Range("E43") '  +2.29%  '

Set-TextValue $ws.Range("E44") '  +5.67%  '

Set-TextValue $ws.Range("E45") '  +6.33%  '

Set-TextValue $ws.Range("D46") '8.08'
Set-TextValue $ws.Range("E46") '  +9.73%  '

Set-TextValue $ws.Range("D47") '16.90'
Set-TextValue $ws.Range("E47") '  +10.18%  '

Set-TextValue $ws.Range("D48") '95.11'
Set-TextValue $ws.Range("E48") '  +5.67%  '

Set-TextValue $ws.Range("D49") '1.426.59'
Set-TextValue $ws.Range("E49") '  +5.70%  '

Set-TextValue $ws.Range("E50") '  +2.95%  '

Set-TextValue $ws.Range("D51") '47.42'
Set-TextValue $ws.Range("E51") '  +4.53%  '

